# Update "tab_Cessazione_unione_civile___convivenze" reference table:
#   - merge the two separate "Decesso ..." rows (id 5 "Decesso unito
#     civilmente" and id 33 "Decesso convivente di fatto") into a single
#     new row "Decesso del convivente/unito civilmente" with id 98
#   - renumber/shift the rows that follow accordingly
#   - widen column B to fit the new (longer) text
#   - leave the active selection on B8, matching the author's last edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: was id 5 "Decesso unito civilmente" -> id 31 "Accordo tra le parti"
$ws.Range("A6").Value = 31
$ws.Range("B6").Value = "Accordo tra le parti"

# Row 7: was id 31 "Accordo tra le parti" -> id 32 "Recesso unilaterale"
$ws.Range("A7").Value = 32
$ws.Range("B7").Value = "Recesso unilaterale"

# Row 8: was id 32 "Recesso unilaterale" -> id 34 "Matrimonio/unione civile"
$ws.Range("A8").Value = 34
$ws.Range("B8").Value = "Matrimonio/unione civile"

# Row 9: was id 33 "Decesso convivente di fatto" -> id 98, new merged text
$ws.Range("A9").Value = 98
$ws.Range("B9").Value = "Decesso del convivente/unito civilmente"

# Row 10 (old id 34 "Matrimonio/unione civile") no longer needed - removed
$ws.Rows.Item(10).Delete()

# Column B needs to be widened to fit the longer description text
$ws.Columns.Item(2).ColumnWidth = 41.33

# Page setup: A4 portrait
$ps = $ws.PageSetup()
$ps.PaperSize = 9
$ps.Orientation = 1

# Match the final selection left by the author
[void]$ws.Range("B8").Select()
